# CU0007 - "Buscar Etiquetas": split the requirement-row description so
# the trailing period becomes its own run (Arial 10pt), matching the
# authored edit ("Agrego interfaces y analisis").

$d = $word.ActiveDocument

# 1) Trim the sentence: drop " a la base de datos" but keep the final period.
$d.Content.Find.Execute(
    "El sistema requiere el filtrado por etiqueta a la base de datos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El sistema requiere el filtrado por etiqueta.", 2)

# 2) Re-find the trimmed sentence and isolate just the trailing "."
$sentence = $d.Content
$sentence.Find.Execute(
    "El sistema requiere el filtrado por etiqueta.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$dot = $d.Range($sentence.End - 1, $sentence.End)

# 3) Give the period explicit Arial/10pt formatting so it becomes its own
#    run, distinct from the sentence that precedes it.
$dot.Font.NameAscii = "Arial"
$dot.Font.NameFarEast = "Arial"
$dot.Font.NameOther = "Arial"
$dot.Font.NameBi = "Arial"
$dot.Font.Size = 10
$dot.Font.Color = -16777216
